$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-22 Monday", "2024-01-23 Tuesday"),
    @("75×84=6300", "79×75=5925"),
    @("76×37=2812", "38×56=2128"),
    @("18×51=918", "75×25=1875"),
    @("52×61=3172", "84×12=1008"),
    @("57×70=3990", "56×61=3416"),
    @("97×63=6111", "17×78=1326"),
    @("28×51=1428", "60×59=3540"),
    @("50×90=4500", "25×60=1500"),
    @("62×58=3596", "57×56=3192"),
    @("46×99=4554", "54×70=3780"),
    @("18×80=1440", "39×35=1365"),
    @("59×99=5841", "82×41=3362"),
    @("25×42=1050", "69×64=4416"),
    @("68×95=6460", "41×20=820"),
    @("17×16=272", "97×27=2619"),
    @("29×45=1305", "46×24=1104"),
    @("39×74=2886", "65×58=3770"),
    @("35×89=3115", "34×94=3196"),
    @("21×34=714", "81×83=6723"),
    @("54×75=4050", "59×49=2891"),
    @("52×69=3588", "11×88=968"),
    @("25×32=800", "13×15=195"),
    @("73×67=4891", "42×81=3402"),
    @("30×39=1170", "40×26=1040"),
    @("94×48=4512", "88×56=4928")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
